$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the row labels to include units ("(seconds)")
$ws.Range("C5").Value = "Cleveland (seconds)"
$ws.Range("C6").Value = "Georgia (seconds)"
$ws.Range("C10").Value = "Cleveland (seconds)"
$ws.Range("C11").Value = "Georgia (seconds)"

# Move the active selection to F13
$ws.Range("F13").Select()
